# obs_max update and clean up
# The deck's "Date Placeholder" fields (datetimeFigureOut) on the slide
# master and every custom (slide) layout were showing the stale cached
# date "12/29/22". Refresh them to "1/11/23" everywhere they appear.

$p = $ppt.ActivePresentation
$oldDate = "12/29/22"
$newDate = "1/11/23"

# --- Slide Master ---
$master = $p.SlideMaster
for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Every Slide Layout ---
$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $lo = $layouts.Item($i)
    for ($j = 1; $j -le $lo.Shapes.Count; $j++) {
        $sh = $lo.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
